$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 816.8
$ws.Range("I38").Value = 526
$ws.Range("J38").Value = 1980
$ws.Range("K38").Value = 1578
$ws.Range("L38").Value = 5940
$ws.Range("M38").Value = -1206
$ws.Range("N38").Value = -6684
$ws.Range("H58").Value = 1949.2222
$ws.Range("I58").Value = 587.6667
$ws.Range("J58").Value = 4672.3335
$ws.Range("K58").Value = 1763.0001
$ws.Range("L58").Value = 14017.0005
$ws.Range("M58").Value = -1613.0001
$ws.Range("N58").Value = -14317.0005
$ws.Range("H100").Value = 3276.25
$ws.Range("I100").Value = 3276.25
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3276.25
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2735.25
$ws.Range("N100").ClearContents()
$ws.Range("H121").Value = 2784.25
$ws.Range("J121").Value = 3662.3333
$ws.Range("L121").Value = 10986.9999
$ws.Range("N121").Value = -14480.9999
$ws.Range("H129").Value = 863.7377300000001
$ws.Range("J129").Value = 887.12726
$ws.Range("L129").Value = 2661.38178
$ws.Range("N129").Value = -12661.38178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4018.1887
$ws.Range("I32").Value = 2673.2827
$ws.Range("K32").Value = 2673.2827
$ws.Range("M32").Value = -2386.2827
$ws.Range("H61").Value = 5843.4
$ws.Range("I61").Value = 4224.5
$ws.Range("J61").Value = 6922.6665
$ws.Range("K61").Value = 4224.5
$ws.Range("L61").Value = 6922.6665
$ws.Range("M61").Value = -4012.5
$ws.Range("N61").Value = -7346.6665
$ws.Range("H74").Value = 1129.7667
$ws.Range("I74").Value = 767.7
$ws.Range("K74").Value = 767.7
$ws.Range("M74").Value = 106.3
$ws.Range("H77").Value = 1129.7667
$ws.Range("I77").Value = 767.7
$ws.Range("K77").Value = 3838.5
$ws.Range("M77").Value = 529.5
$ws.Range("H109").Value = 47916.668
$ws.Range("J109").Value = 47916.668
$ws.Range("L109").Value = 47916.668
$ws.Range("N109").Value = -50690.668
$ws.Range("H132").Value = 1564.5106
$ws.Range("I132").Value = 1030.1177
$ws.Range("K132").Value = 3090.3531
$ws.Range("M132").Value = -560.3531000000003
$ws.Range("H136").Value = 5843.4
$ws.Range("I136").Value = 4224.5
$ws.Range("J136").Value = 6922.6665
$ws.Range("K136").Value = 12673.5
$ws.Range("L136").Value = 20767.9995
$ws.Range("M136").Value = -10123.5
$ws.Range("N136").Value = -25867.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 136625.8
$ws.Range("I86").Value = 6137.6
$ws.Range("J86").Value = 201869.9
$ws.Range("K86").Value = 6137.6
$ws.Range("L86").Value = 201869.9
$ws.Range("M86").Value = -5014.6
$ws.Range("N86").Value = -204115.9
$ws.Range("H89").Value = 136625.8
$ws.Range("I89").Value = 6137.6
$ws.Range("J89").Value = 201869.9
$ws.Range("K89").Value = 30688
$ws.Range("L89").Value = 1009349.5
$ws.Range("M89").Value = -25072
$ws.Range("N89").Value = -1020581.5
$ws.Range("H105").Value = 2500
$ws.Range("I105").Value = 2500
$ws.Range("K105").Value = 2500
$ws.Range("M105").Value = -753
$ws.Range("H107").Value = 1179.5
$ws.Range("I107").Value = 1195.7273
$ws.Range("J107").Value = 1001
$ws.Range("K107").Value = 1195.7273
$ws.Range("L107").Value = 1001
$ws.Range("M107").Value = 724.2727
$ws.Range("N107").Value = -4841
$ws.Range("H108").Value = 94995
$ws.Range("J108").Value = 94995
$ws.Range("L108").Value = 94995
$ws.Range("N108").Value = -102675
$ws.Range("H134").Value = 9545.333000000001
$ws.Range("I134").Value = 14365.182
$ws.Range("K134").Value = 43095.546
$ws.Range("M134").Value = -40560.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2551.8096
$ws.Range("I31").Value = 2490.2222
$ws.Range("K31").Value = 2490.2222
$ws.Range("M31").Value = -2195.2222
$ws.Range("H34").Value = 2551.8096
$ws.Range("I34").Value = 2490.2222
$ws.Range("K34").Value = 2490.2222
$ws.Range("M34").Value = -2288.2222
$ws.Range("H54").Value = 16000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H58").Value = 1554336.6
$ws.Range("I58").Value = 2289400
$ws.Range("K58").Value = 2289400
$ws.Range("M58").Value = -2289197
$ws.Range("H122").Value = 3437
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3437
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 10311
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -15211
$ws.Range("H134").Value = 1204.5103
$ws.Range("I134").Value = 1221.6279
$ws.Range("J134").Value = 1081.8334
$ws.Range("K134").Value = 3664.8837
$ws.Range("L134").Value = 3245.5002
$ws.Range("M134").Value = -1129.8837
$ws.Range("N134").Value = -8315.5002
$ws.Range("H136").Value = 1554336.6
$ws.Range("I136").Value = 2289400
$ws.Range("K136").Value = 6868200
$ws.Range("M136").Value = -6865650
$ws.Range("H141").Value = 71552.8
$ws.Range("J141").Value = 71552.8
$ws.Range("L141").Value = 71552.8
$ws.Range("N141").Value = -81912.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 602.7692
$ws.Range("I68").Value = 564.6667
$ws.Range("J68").Value = 635.4286
$ws.Range("K68").Value = 1694.0001
$ws.Range("L68").Value = 1906.2858
$ws.Range("M68").Value = -883.0001
$ws.Range("N68").Value = -3528.2858
$ws.Range("H71").Value = 602.7692
$ws.Range("I71").Value = 564.6667
$ws.Range("J71").Value = 635.4286
$ws.Range("K71").Value = 5082.0003
$ws.Range("L71").Value = 5718.8574
$ws.Range("M71").Value = -1026.0003
$ws.Range("N71").Value = -13830.8574
$ws.Range("H122").Value = 884.7857
$ws.Range("I122").Value = 453.25
$ws.Range("J122").Value = 1057.4
$ws.Range("K122").Value = 4079.25
$ws.Range("L122").Value = 9516.6
$ws.Range("M122").Value = -1629.25
$ws.Range("N122").Value = -14416.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H102").Value = 1960.7858
$ws.Range("I102").Value = 1996.2693
$ws.Range("K102").Value = 1996.2693
$ws.Range("M102").Value = -374.2692999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2569.15
$ws.Range("J7").Value = 9998.5
$ws.Range("L7").Value = 9998.5
$ws.Range("N7").Value = -10222.5
$ws.Range("H22").Value = 2838.7778
$ws.Range("I22").Value = 5116.6665
$ws.Range("J22").Value = 1699.8334
$ws.Range("K22").Value = 5116.6665
$ws.Range("L22").Value = 1699.8334
$ws.Range("M22").Value = -4821.6665
$ws.Range("N22").Value = -2289.8334
$ws.Range("H27").Value = 2838.7778
$ws.Range("I27").Value = 5116.6665
$ws.Range("J27").Value = 1699.8334
$ws.Range("K27").Value = 5116.6665
$ws.Range("L27").Value = 1699.8334
$ws.Range("M27").Value = -5009.6665
$ws.Range("N27").Value = -1913.8334
$ws.Range("H40").Value = 3023.3572
$ws.Range("I40").Value = 2312.0908
$ws.Range("J40").Value = 5631.3335
$ws.Range("K40").Value = 2312.0908
$ws.Range("L40").Value = 5631.3335
$ws.Range("M40").Value = -2176.0908
$ws.Range("N40").Value = -5903.3335
$ws.Range("H43").Value = 10808.4
$ws.Range("J43").Value = 10808.4
$ws.Range("L43").Value = 10808.4
$ws.Range("N43").Value = -11194.4
$ws.Range("H126").Value = 2569.15
$ws.Range("J126").Value = 9998.5
$ws.Range("L126").Value = 29995.5
$ws.Range("N126").Value = -34935.5
$ws.Range("H130").Value = 25000
$ws.Range("J130").Value = 25000
$ws.Range("L130").Value = 25000
$ws.Range("N130").Value = -35040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1116.6316
$ws.Range("J14").Value = 1116.6316
$ws.Range("L14").Value = 1116.6316
$ws.Range("N14").Value = -1452.6316
$ws.Range("H81").Value = 2685
$ws.Range("J81").Value = 2600
$ws.Range("L81").Value = 5200
$ws.Range("N81").Value = -7322
$ws.Range("H84").Value = 2685
$ws.Range("J84").Value = 2600
$ws.Range("L84").Value = 26000
$ws.Range("N84").Value = -36608
$ws.Range("H113").Value = 505.1905
$ws.Range("I113").Value = 321
$ws.Range("K113").Value = 963
$ws.Range("M113").Value = 1207
$ws.Range("H123").Value = 47599.75
$ws.Range("J123").Value = 47599.75
$ws.Range("L123").Value = 47599.75
$ws.Range("N123").Value = -57399.75
$ws.Range("H126").Value = 2891.8572
$ws.Range("I126").Value = 2032.25
$ws.Range("K126").Value = 6096.75
$ws.Range("M126").Value = -3626.75
